$wb = $excel.ActiveWorkbook

# --- Fill in new test-id values in column A (order matters: matches shared-string allocation order) ---
$ws = $wb.Worksheets.Item("getSensorByDeviceId")
$ws.Range("A2").Value = "iems-api-service-sp5-6"
$ws.Range("A5").Value = "iems-api-service-sp5-7"
$ws.Range("A6").Value = "iems-api-service-sp5-8"
$ws = $wb.Worksheets.Item("getSensorDataBySensorId")
$ws.Range("A2").Value = "iems-api-service-sp5-9"
$ws.Range("A4").Value = "iems-api-service-sp5-12"
$ws.Range("A3").Value = "iems-api-service-sp5-13"
$ws = $wb.Worksheets.Item("getSensorDataByDeviceId")
$ws.Range("A2").Value = "iems-api-service-sp5-17"
$ws.Range("A4").Value = "iems-api-service-sp5-20"
$ws.Range("A3").Value = "iems-api-service-sp5-21"
$ws = $wb.Worksheets.Item("subscriptionsByDeviceId")
$ws.Range("A2").Value = "iems-api-service-sp5-24"
$ws = $wb.Worksheets.Item("deleteSubscriptions")
$ws.Range("A2").Value = "iems-api-service-sp5-26"
$ws.Range("A3").Value = "iems-api-service-sp5-27"
$ws = $wb.Worksheets.Item("subscriptionsBySensorId")
$ws.Range("A2").Value = "iems-api-service-sp5-28"
$ws.Range("A4").Value = "iems-api-service-sp5-29"
$ws = $wb.Worksheets.Item("subscriptionsWithKPIByDeviceId")
$ws.Range("A2").Value = "iems-api-service-kpi-2"
$ws = $wb.Worksheets.Item("getTopSensorDataByDeviceId")
$ws.Range("A2").Value = "iems-api-service-kpi-6"
$ws.Range("A6").Value = "iems-api-service-kpi-7"
$ws.Range("A5").Value = "iems-api-service-kpi-8"
$ws.Range("A4").Value = "iems-api-service-kpi-11"
$ws = $wb.Worksheets.Item("getTopKPIDataByDeviceId")
$ws.Range("A2").Value = "iems-api-service-kpi-12"
$ws.Range("A6").Value = "iems-api-service-kpi-13"
$ws.Range("A5").Value = "iems-api-service-kpi-14"
$ws.Range("A4").Value = "iems-api-service-kpi-17"
$ws = $wb.Worksheets.Item("getKpiDataByDeviceId")
$ws.Range("A2").Value = "iems-api-service-kpi-18"
$ws.Range("A4").Value = "iems-api-service-kpi-21"
$ws.Range("A3").Value = "iems-api-service-kpi-22"
$ws.Range("A6").Value = "iems-api-service-kpi-24"
$ws = $wb.Worksheets.Item("getDeviceInfoByID")
$ws.Range("A2").Value = "iems-api-service-kpi-26"
$ws.Range("A5").Value = "iems-api-service-kpi-27"
$ws.Range("A3").Value = "iems-api-service-kpi-29-var1"
$ws.Range("A4").Value = "iems-api-service-kpi-29-var2"
$ws = $wb.Worksheets.Item("getTopSensorDataByDeviceId")
$ws.Range("A3").Value = "iems-api-service-kpi-35"

# --- Restore per-sheet selections (active cell) ---
$wb.Worksheets.Item("getDevicesByType").Range("A2").Select()
$wb.Worksheets.Item("getDeviceInfoByID").Range("A6").Select()
$wb.Worksheets.Item("getSensorDataBySensorId").Range("B6").Select()
$wb.Worksheets.Item("getSensorDataByDeviceId").Range("B5").Select()
$wb.Worksheets.Item("getTopSensorDataByDeviceId").Range("A3").Select()
$wb.Worksheets.Item("getKpiDataByDeviceId").Range("D21").Select()
$wb.Worksheets.Item("getTopKPIDataByDeviceId").Range("A4").Select()
$wb.Worksheets.Item("subscriptionsBySensorId").Range("B5").Select()
$wb.Worksheets.Item("subscriptionsByDeviceId").Range("A3").Select()
$wb.Worksheets.Item("subscriptionsWithKPIByDeviceId").Range("A5").Select()
$wb.Worksheets.Item("deleteSubscriptions").Range("A3").Select()

# --- Activate getSensorByDeviceId last so it becomes the active/selected tab (matches activeTab=2) ---
$ws3 = $wb.Worksheets.Item("getSensorByDeviceId")
$ws3.Activate()
$ws3.Range("A4").Select()

